# Insert a new data row at row 738 (shifts existing rows 738:839 down to 739:840)
# and populate it with a new daily price record, matching the author's commit
# "Fruta / hortaliza, semanal" (a new weekly data point was added to the series).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("738:738").Insert()

$ws.Cells.Item(738, 1).Value2  = 11
$ws.Cells.Item(738, 2).Value2  = "Vega Monumental Concepción"
$ws.Cells.Item(738, 3).Value2  = "Bíobío"
$ws.Cells.Item(738, 4).Value2  = 45142
$ws.Cells.Item(738, 5).Value2  = 8
$ws.Cells.Item(738, 6).Value2  = "Fruta"
$ws.Cells.Item(738, 7).Value2  = 100102
$ws.Cells.Item(738, 8).Value2  = "Cítricos"
$ws.Cells.Item(738, 9).Value2  = 100102003
$ws.Cells.Item(738, 10).Value2 = "Limón"
$ws.Cells.Item(738, 11).Value2 = "Sin especificar"
$ws.Cells.Item(738, 12).Value2 = "1a amarillo"
$ws.Cells.Item(738, 13).Value2 = 500
$ws.Cells.Item(738, 14).Value2 = 6000
$ws.Cells.Item(738, 15).Value2 = 6500
$ws.Cells.Item(738, 16).Value2 = 6300
$ws.Cells.Item(738, 17).Value2 = "$/malla 18 kilos"
$ws.Cells.Item(738, 18).Value2 = "Provincia de Limarí"
$ws.Cells.Item(738, 19).Value2 = 350
$ws.Cells.Item(738, 20).Value2 = 18
